# Add a "Correct_answer" column (D) to the stimulus table.
#   D1            -> header "Correct_answer"
#   D2:D29  (Purple rows) -> "l"
#   D30:D101 (Blue rows)  -> "s"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Correct_answer"
$ws.Range("D2:D29").Value = "l"
$ws.Range("D30:D101").Value = "s"

# Restore the view state as closely as possible: the author's selection
# ends up on D30:D101 with the viewport scrolled so row 84 is visible.
$ws.Range("D30:D101").Select()
